$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G (row 4 and 5) value changes
$ws.Cells.Item(4, 7).Value = 5
$ws.Cells.Item(5, 7).Value = 4

# Column A toggled values (row -> new value)
$aChanges = @(
    @(2, 0),
    @(5, 0),
    @(6, 1),
    @(12, 0),
    @(18, 1),
    @(23, 1),
    @(24, 1),
    @(30, 0),
    @(31, 0),
    @(33, 1),
    @(37, 1),
    @(38, 0),
    @(44, 1),
    @(45, 0),
    @(46, 1),
    @(51, 0),
    @(53, 0),
    @(54, 0),
    @(59, 0),
    @(61, 1),
    @(65, 1),
    @(66, 1),
    @(73, 1),
    @(78, 0),
    @(79, 0),
    @(80, 1),
    @(84, 0),
    @(85, 0),
    @(86, 1),
    @(90, 0),
    @(92, 1),
    @(93, 0),
    @(97, 1),
    @(104, 1),
    @(105, 0),
    @(106, 1),
    @(111, 0),
    @(113, 1),
    @(116, 0),
    @(119, 1),
    @(132, 0),
    @(133, 1),
    @(134, 0),
    @(137, 1),
    @(145, 1),
    @(152, 1),
    @(155, 1),
    @(159, 1),
    @(160, 1),
    @(162, 0),
    @(166, 0),
    @(168, 0),
    @(169, 1),
    @(172, 0),
    @(173, 0),
    @(174, 0),
    @(175, 1),
    @(176, 0),
    @(178, 0),
    @(184, 1),
    @(186, 0),
    @(187, 1),
    @(193, 0),
    @(194, 1),
    @(198, 0),
    @(200, 1),
    @(205, 1),
    @(206, 1),
    @(207, 0),
    @(209, 0),
    @(212, 0),
    @(213, 1),
    @(215, 1),
    @(219, 0),
    @(220, 1),
    @(222, 0),
    @(227, 0),
    @(228, 1),
    @(230, 1),
    @(233, 1),
    @(236, 0),
    @(237, 0),
    @(242, 1),
    @(243, 0),
    @(247, 1),
    @(251, 0),
    @(255, 0),
    @(260, 1),
    @(262, 1),
    @(266, 1),
    @(267, 0),
    @(274, 1),
    @(275, 0),
    @(281, 0),
    @(286, 1),
    @(291, 1),
    @(294, 0),
    @(295, 0),
    @(305, 0),
    @(307, 1),
    @(311, 0),
    @(316, 1),
    @(320, 0),
    @(329, 1),
    @(334, 0),
    @(335, 0),
    @(336, 0),
    @(337, 1),
    @(341, 1),
    @(344, 0),
    @(348, 1),
    @(349, 1),
    @(352, 1),
    @(353, 1),
    @(355, 1),
    @(359, 0),
    @(361, 1),
    @(369, 0),
    @(374, 0),
    @(375, 0),
    @(382, 1),
    @(383, 0),
    @(390, 0),
    @(396, 0),
    @(397, 1),
    @(401, 1),
    @(402, 1),
    @(404, 1),
    @(408, 0),
    @(409, 1),
    @(415, 0),
    @(416, 0),
    @(417, 0),
    @(418, 1),
    @(422, 0),
    @(423, 0),
    @(428, 1),
    @(429, 1),
    @(434, 1),
    @(436, 0),
    @(437, 0),
    @(441, 0),
    @(443, 1),
    @(449, 0),
    @(450, 1),
    @(455, 1),
    @(463, 0),
    @(470, 1),
    @(483, 0),
    @(484, 1),
    @(486, 1),
    @(491, 0),
    @(502, 0),
    @(504, 0),
    @(505, 0),
    @(511, 1),
    @(515, 1),
    @(518, 1),
    @(525, 0),
    @(526, 1),
    @(527, 1),
    @(529, 1),
    @(530, 0),
    @(535, 0),
    @(536, 0),
    @(544, 0),
    @(545, 0),
    @(555, 1),
    @(556, 1),
    @(557, 1),
    @(569, 1),
    @(572, 1),
    @(579, 1),
    @(582, 0),
    @(586, 0),
    @(590, 0),
    @(597, 0),
    @(605, 1),
    @(608, 1),
    @(609, 1),
    @(626, 0),
    @(629, 0),
    @(639, 0),
    @(643, 0),
    @(645, 0),
    @(650, 1),
    @(652, 1),
    @(655, 1),
    @(656, 1),
    @(658, 1),
    @(662, 0),
    @(663, 0),
    @(664, 1),
    @(665, 0),
    @(669, 0),
    @(670, 1),
    @(676, 1),
    @(677, 0),
    @(679, 0),
    @(680, 1),
    @(690, 0),
    @(693, 0),
    @(698, 1),
    @(700, 1),
    @(702, 0),
    @(704, 0),
    @(705, 1),
    @(709, 0),
    @(712, 0),
    @(719, 1),
    @(720, 1),
    @(733, 1),
    @(745, 1),
    @(746, 1),
    @(754, 0),
    @(757, 1),
    @(761, 1),
    @(766, 0),
    @(767, 0),
    @(768, 0),
    @(772, 1),
    @(779, 1),
    @(785, 1),
    @(800, 0),
    @(805, 0),
    @(806, 0),
    @(815, 0),
    @(816, 0),
    @(819, 0),
    @(825, 0),
    @(827, 0),
    @(832, 1),
    @(834, 1),
    @(835, 1),
    @(836, 1),
    @(837, 1),
    @(839, 0),
    @(840, 1),
    @(844, 1),
    @(845, 0),
    @(852, 1),
    @(854, 1),
    @(855, 0),
    @(866, 0),
    @(867, 0),
    @(869, 0),
    @(874, 1),
    @(876, 1)
)

foreach ($pair in $aChanges) {
    $ws.Cells.Item($pair[0], 1).Value = $pair[1]
}

# Update selection to J3
$ws.Range("J3").Select() | Out-Null
